# UnitType.xlsx: add a "unitTypeId" column right after "projectId".
#
# Before: projectId | name    | sellingPrice | available | total
# After:  projectId | unitTypeId | name | sellingPrice | available | total
#
# The new column gets a simple 0-based identifier per row (0, 1, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B, shifting name/sellingPrice/available/total right.
$ws.Columns("B:B").Insert()

# Header for the new column.
$ws.Range("B1").Value = "unitTypeId"

# Populate the new unitTypeId values for the existing data rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 2
}
